# Auto-generated Excel COM-interop script applying the Tonberry_Profits value updates.
# For each sheet, update the H..N "market price / profit" columns for the affected Leve rows.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 6: H6,I6,J6,K6,L6,M6,N6
$ws.Range("H6").Value = 600.5714
$ws.Range("I6").Value = 50.75
$ws.Range("J6").Value = 1333.6666
$ws.Range("K6").Value = 152.25
$ws.Range("L6").Value = 4000.9998
$ws.Range("M6").Value = -40.25
$ws.Range("N6").Value = -4224.9998
# Row 28: H28,J28,L28,N28
$ws.Range("H28").Value = 498.5
$ws.Range("J28").Value = 998.5
$ws.Range("L28").Value = 998.5
$ws.Range("N28").Value = -1968.5
# Row 33: H33,J33,L33,N33
$ws.Range("H33").Value = 259.4091
$ws.Range("J33").Value = 116.9
$ws.Range("L33").Value = 116.9
$ws.Range("N33").Value = -574.9
# Row 42: H42,I42,J42,K42,L42,M42,N42
$ws.Range("H42").Value = 193.625
$ws.Range("I42").Value = 116.666664
$ws.Range("J42").Value = 239.8
$ws.Range("K42").Value = 349.999992
$ws.Range("L42").Value = 719.4000000000001
$ws.Range("M42").Value = -119.999992
$ws.Range("N42").Value = -1179.4
# Row 43: H43,I43,J43,K43,L43,M43,N43
$ws.Range("H43").Value = 891.2
$ws.Range("I43").Value = 724.5
$ws.Range("J43").Value = 962.6429000000001
$ws.Range("K43").Value = 724.5
$ws.Range("L43").Value = 962.6429000000001
$ws.Range("M43").Value = -655.5
$ws.Range("N43").Value = -1100.6429
# Row 62: H62,I62,J62,K62,L62,M62,N62
$ws.Range("H62").Value = 1672
$ws.Range("I62").Value = 1650.6666
$ws.Range("J62").Value = 1800
$ws.Range("K62").Value = 1650.6666
$ws.Range("L62").Value = 1800
$ws.Range("M62").Value = -1026.6666
$ws.Range("N62").Value = -3048
# Row 65: H65,I65,J65,K65,L65,M65,N65
$ws.Range("H65").Value = 1672
$ws.Range("I65").Value = 1650.6666
$ws.Range("J65").Value = 1800
$ws.Range("K65").Value = 8253.333000000001
$ws.Range("L65").Value = 9000
$ws.Range("M65").Value = -5133.333000000001
$ws.Range("N65").Value = -15240
# Row 92: H92,I92,K92,M92
$ws.Range("H92").Value = 1368490.5
$ws.Range("I92").Value = 2052402.9
$ws.Range("K92").Value = 2052402.9
$ws.Range("M92").Value = -2051154.9
# Row 129: H129,J129,L129,N129
$ws.Range("H129").Value = 863.5685999999999
$ws.Range("J129").Value = 875.59576
$ws.Range("L129").Value = 2626.78728
$ws.Range("N129").Value = -12626.78728
# Row 138: H138,I138,J138,K138,L138,M138,N138
$ws.Range("H138").Value = 3011.132
$ws.Range("I138").Value = 3004.2173
$ws.Range("J138").Value = 3016.4333
$ws.Range("K138").Value = 9012.651899999999
$ws.Range("L138").Value = 9049.2999
$ws.Range("M138").Value = -3872.651899999999
$ws.Range("N138").Value = -19329.2999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2: H2,I2,J2,K2,L2,M2,N2
$ws.Range("H2").Value = 5555555
$ws.Range("I2").Value = 5555555
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 5555555
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -5555442
$ws.Range("N2").ClearContents()
# Row 32: H32,I32,J32,K32,L32,M32,N32
$ws.Range("H32").Value = 3537.4695
$ws.Range("I32").Value = 3197.9302
$ws.Range("J32").Value = 5970.8335
$ws.Range("K32").Value = 3197.9302
$ws.Range("L32").Value = 5970.8335
$ws.Range("M32").Value = -2910.9302
$ws.Range("N32").Value = -6544.8335
# Row 45: H45,I45,K45,M45
$ws.Range("H45").Value = 1669.6875
$ws.Range("I45").Value = 1440.375
$ws.Range("K45").Value = 1440.375
$ws.Range("M45").Value = -1063.375
# Row 61: H61,I61,J61,K61,L61,M61,N61
$ws.Range("H61").Value = 5857.143
$ws.Range("I61").Value = 1250
$ws.Range("J61").Value = 12000
$ws.Range("K61").Value = 1250
$ws.Range("L61").Value = 12000
$ws.Range("M61").Value = -1038
$ws.Range("N61").Value = -12424
# Row 74: H74,I74,K74,M74
$ws.Range("H74").Value = 1570.3158
$ws.Range("I74").Value = 1404.5454
$ws.Range("K74").Value = 1404.5454
$ws.Range("M74").Value = -530.5454
# Row 77: H77,I77,K77,M77
$ws.Range("H77").Value = 1570.3158
$ws.Range("I77").Value = 1404.5454
$ws.Range("K77").Value = 7022.727
$ws.Range("M77").Value = -2654.727
# Row 97: H97,I97,K97,M97
$ws.Range("H97").Value = 1284.6666
$ws.Range("I97").Value = 1284.6666
$ws.Range("K97").Value = 1284.6666
$ws.Range("M97").Value = -788.6666
# Row 110: H110,I110,K110,M110
$ws.Range("H110").Value = 2320.3704
$ws.Range("I110").Value = 1670.05
$ws.Range("K110").Value = 1670.05
$ws.Range("M110").Value = 374.95
# Row 116: H116,I116,J116,K116,L116,M116,N116
$ws.Range("H116").Value = 5555555
$ws.Range("I116").Value = 5555555
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 5555555
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -5553261
$ws.Range("N116").ClearContents()
# Row 132: H132,I132,J132,K132,L132,M132,N132
$ws.Range("H132").Value = 1720.2632
$ws.Range("I132").Value = 1068.6333
$ws.Range("J132").Value = 4163.875
$ws.Range("K132").Value = 3205.8999
$ws.Range("L132").Value = 12491.625
$ws.Range("M132").Value = -675.8998999999999
$ws.Range("N132").Value = -17551.625
# Row 136: H136,I136,J136,K136,L136,M136,N136
$ws.Range("H136").Value = 5857.143
$ws.Range("I136").Value = 1250
$ws.Range("J136").Value = 12000
$ws.Range("K136").Value = 3750
$ws.Range("L136").Value = 36000
$ws.Range("M136").Value = -1200
$ws.Range("N136").Value = -41100

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3: H3,I3,J3,K3,L3,M3,N3
$ws.Range("H3").Value = 5555555
$ws.Range("I3").Value = 5555555
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 5555555
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -5555441
$ws.Range("N3").ClearContents()
# Row 20: H20,I20,K20,M20
$ws.Range("H20").Value = 2176.8276
$ws.Range("I20").Value = 2117.375
$ws.Range("K20").Value = 2117.375
$ws.Range("M20").Value = -1870.375
# Row 94: H94,I94,K94,M94
$ws.Range("H94").Value = 791.9231
$ws.Range("I94").Value = 626.8182
$ws.Range("K94").Value = 626.8182
$ws.Range("M94").Value = -175.8182
# Row 99: H99,I99,K99,M99
$ws.Range("H99").Value = 1901.3889
$ws.Range("I99").Value = 1680.5385
$ws.Range("K99").Value = 1680.5385
$ws.Range("M99").Value = -182.5385000000001
# Row 107: H107,I107,J107,K107,L107,M107,N107
$ws.Range("H107").Value = 1814.625
$ws.Range("I107").Value = 1545.1538
$ws.Range("J107").Value = 2133.0908
$ws.Range("K107").Value = 1545.1538
$ws.Range("L107").Value = 2133.0908
$ws.Range("M107").Value = 374.8462
$ws.Range("N107").Value = -5973.0908
# Row 134: H134,I134,K134,M134
$ws.Range("H134").Value = 6227.212
$ws.Range("I134").Value = 6984.615
$ws.Range("K134").Value = 20953.845
$ws.Range("M134").Value = -18418.845

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31: H31,I31,K31,M31
$ws.Range("H31").Value = 2426.625
$ws.Range("I31").Value = 2363.3
$ws.Range("K31").Value = 2363.3
$ws.Range("M31").Value = -2068.3
# Row 34: H34,I34,K34,M34
$ws.Range("H34").Value = 2426.625
$ws.Range("I34").Value = 2363.3
$ws.Range("K34").Value = 2363.3
$ws.Range("M34").Value = -2161.3
# Row 132: H132,I132,J132,K132,L132,M132,N132
$ws.Range("H132").Value = 2023.04
$ws.Range("I132").Value = 1089.75
$ws.Range("J132").Value = 5756.2
$ws.Range("K132").Value = 3269.25
$ws.Range("L132").Value = 17268.6
$ws.Range("M132").Value = -739.25
$ws.Range("N132").Value = -22328.6
# Row 134: H134,I134,K134,M134
$ws.Range("H134").Value = 2368.9583
$ws.Range("I134").Value = 2097.15
$ws.Range("K134").Value = 6291.450000000001
$ws.Range("M134").Value = -3756.450000000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 7: H7,J7,L7,N7
$ws.Range("H7").Value = 306.14285
$ws.Range("J7").Value = 371.55554
$ws.Range("L7").Value = 1114.66662
$ws.Range("N7").Value = -1338.66662
# Row 15: H15,I15,J15,K15,L15,M15,N15
$ws.Range("H15").Value = 50
$ws.Range("I15").Value = 50
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 150
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -10
$ws.Range("N15").ClearContents()
# Row 121: H121,I121,J121,K121,L121,M121,N121
$ws.Range("H121").Value = 595.8823
$ws.Range("I121").Value = 499.81818
$ws.Range("J121").Value = 772
$ws.Range("K121").Value = 1499.45454
$ws.Range("L121").Value = 2316
$ws.Range("M121").Value = -189.45454
$ws.Range("N121").Value = -4936

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80: H80,I80,J80,K80,L80,M80,N80
$ws.Range("H80").Value = 2437.5
$ws.Range("I80").Value = 2340
$ws.Range("J80").Value = 2600
$ws.Range("K80").Value = 2340
$ws.Range("L80").Value = 2600
$ws.Range("M80").Value = -1342
$ws.Range("N80").Value = -4596
# Row 83: H83,I83,J83,K83,L83,M83,N83
$ws.Range("H83").Value = 2437.5
$ws.Range("I83").Value = 2340
$ws.Range("J83").Value = 2600
$ws.Range("K83").Value = 11700
$ws.Range("L83").Value = 13000
$ws.Range("M83").Value = -6708
$ws.Range("N83").Value = -22984
# Row 102: H102,I102,J102,K102,L102,M102,N102
$ws.Range("H102").Value = 2247.158
$ws.Range("I102").Value = 2274.6775
$ws.Range("J102").Value = 2125.2856
$ws.Range("K102").Value = 2274.6775
$ws.Range("L102").Value = 2125.2856
$ws.Range("M102").Value = -652.6774999999998
$ws.Range("N102").Value = -5369.2856
# Row 132: H132,I132,J132,K132,L132,M132,N132
$ws.Range("H132").Value = 1541186
$ws.Range("I132").Value = 2264388.2
$ws.Range("J132").Value = 4381.125
$ws.Range("K132").Value = 6793164.600000001
$ws.Range("L132").Value = 13143.375
$ws.Range("M132").Value = -6790634.600000001
$ws.Range("N132").Value = -18203.375

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7: H7,I7,J7,K7,L7,M7,N7
$ws.Range("H7").Value = 3186.2666
$ws.Range("I7").Value = 2199.3333
$ws.Range("J7").Value = 4666.6665
$ws.Range("K7").Value = 2199.3333
$ws.Range("L7").Value = 4666.6665
$ws.Range("M7").Value = -2087.3333
$ws.Range("N7").Value = -4890.6665
# Row 93: H93,I93,J93,K93,L93,M93,N93
$ws.Range("H93").Value = 1600
$ws.Range("I93").Value = 1000
$ws.Range("J93").Value = 2500
$ws.Range("K93").Value = 1000
$ws.Range("L93").Value = 2500
$ws.Range("M93").Value = 248
$ws.Range("N93").Value = -4996
# Row 95: H95,J95,L95,N95
$ws.Range("H95").Value = 49999
$ws.Range("J95").Value = 49999
$ws.Range("L95").Value = 49999
$ws.Range("N95").Value = -55491
# Row 126: H126,I126,J126,K126,L126,M126,N126
$ws.Range("H126").Value = 3186.2666
$ws.Range("I126").Value = 2199.3333
$ws.Range("J126").Value = 4666.6665
$ws.Range("K126").Value = 6597.999899999999
$ws.Range("L126").Value = 13999.9995
$ws.Range("M126").Value = -4127.999899999999
$ws.Range("N126").Value = -18939.9995
# Row 132: H132,I132,J132,K132,L132,M132,N132
$ws.Range("H132").Value = 1707.8857
$ws.Range("I132").Value = 1036.36
$ws.Range("J132").Value = 3386.7
$ws.Range("K132").Value = 3109.08
$ws.Range("L132").Value = 10160.1
$ws.Range("M132").Value = -579.0799999999999
$ws.Range("N132").Value = -15220.1
# Row 136: H136,I136,K136,M136
$ws.Range("H136").Value = 2537.257
$ws.Range("I136").Value = 1375.1111
$ws.Range("K136").Value = 4125.3333
$ws.Range("M136").Value = -1575.3333

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 122: H122,I122,K122,M122
$ws.Range("H122").Value = 40388
$ws.Range("I122").Value = 57225.715
$ws.Range("K122").Value = 171677.145
$ws.Range("M122").Value = -169227.145
# Row 136: H136,I136,K136,M136
$ws.Range("H136").Value = 18520148
$ws.Range("I136").Value = 27778988
$ws.Range("K136").Value = 83336964
$ws.Range("M136").Value = -83334414

